# Applies the "added very basic date/week filter and retained browsing dirs"
# commit: appends a couple of journal sentences to two existing paragraphs
# and adds a brand-new trailing paragraph.
#
# Note: this engine coalesces freshly-appended text into the *previous*
# run whenever a paragraph currently holds exactly one run (it only keeps
# genuinely separate <w:r> elements once 2+ runs already exist). To end up
# with the same run layout as the target document we insert the full
# text blob first, then "saw" it apart at the desired run boundaries using
# a harmless Bold-on/Bold-off round trip (which forces a fresh run there,
# and because formatting reverts to identical state on both sides no
# stray formatting is left behind once 2+ runs exist).

function Split-RunBoundary($doc, $pos, $paraEndPos) {
    $r = $doc.Range($pos, $paraEndPos)
    $r.Bold = 1
    $r.Bold = 0
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Paragraph: "...the browse folder need to be done first or else it wont
# be populated!" -> append two new sentences as two separate runs.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("the browse folder need to be done first or else it wont be populated!", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(" THIS ALSO WORKS PERFECTLY WITH FILTERS WOOOO")
$rng.LanguageID = "en-GB"
$rng.Collapse(0)
$rng.InsertAfter(" because it uses current layout of folder populated including the filters NICE!")
$rng.LanguageID = "en-GB"

# ---------------------------------------------------------------------
# Paragraph: "Now adding dropdown choice for weeks folders etc." ->
# append two more sentences (kept as separate runs from the original
# text and from each other).
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Now adding dropdown choice for weeks folders etc.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Collapse(0)
$para2Start = $rng2.Start

$part2a = " Ok bare basic done, ideally should be dynamic and have the "
$part2b = "disabled according to folder directories but claude reached limits rn so I" + [char]0x2019 + "ll just continue with trying to integrate this with VR first instead"

$rng2.InsertAfter($part2a + $part2b)
$rng2.LanguageID = "en-GB"
$para2End = $rng2.End

# Boundary between part2a and part2b (rightmost split first so earlier
# offsets stay valid).
$boundaryAB = $para2Start + $part2a.Length
Split-RunBoundary $d $boundaryAB $para2End

# Boundary between the original "Now adding ... etc." run and part2a.
Split-RunBoundary $d $para2Start $boundaryAB

# ---------------------------------------------------------------------
# New trailing paragraph: "Also, added the browse dir retain
# functionality." split into three runs: "A" / "lso" / ", added the
# browse dir retain functionality."
# ---------------------------------------------------------------------
# Re-use the tracked end of the (now extended) "Now adding ..." paragraph
# instead of re-running Find, since Find would only locate the original
# short sentence and land *before* the text we just appended above.
$rng3 = $d.Range($para2End, $para2End)
$rng3.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last.Range
$newParaStart = $newPara.Start
$fullText = "Also, added the browse dir retain functionality."
$newPara.InsertAfter($fullText)
$newPara.LanguageID = "en-GB"
$newParaEnd = $d.Paragraphs.Last.Range.End

# Split boundaries, rightmost first: "A" | "lso" | ", added..."
$boundary2 = $newParaStart + 4   # right after "Also"
Split-RunBoundary $d $boundary2 $newParaEnd

$boundary1 = $newParaStart + 1   # right after "A"
Split-RunBoundary $d $boundary1 $boundary2

Write-Output "done"
